$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "test"
$ws.Range("B10").Value = "livehta_886_data"
$ws.Range("C10").Value = "\Testdata\Non_Oncology\DataFiles\LIVEHTA_886\livehta_886_Data.xlsx"

$ws.Range("A10:C10").Select()
